# Menu hover y ventanas modales en el listado de activos
# Append new "activo" log rows (126-132) to the registro sheet, mirroring the
# existing B/C "estado" columns (stored as text, e.g. "1"/"2") and the D/E
# timestamp + F note columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rows = @(
    @{ a = 125; b = "1"; c = "1"; d = "Tue Sep 22 22:28:06 CEST 2020"; e = "Tue Sep 22 22:53:52 CEST 2020"; f = "activo 1" },
    @{ a = 126; b = "2"; c = "2"; d = "Tue Sep 22 22:54:14 CEST 2020"; e = "Tue Sep 22 22:54:38 CEST 2020"; f = "de espera activo 2" },
    @{ a = 127; b = "1"; c = "1"; d = "Tue Sep 22 22:53:52 CEST 2020"; e = "Tue Sep 22 22:56:15 CEST 2020"; f = "null" },
    @{ a = 128; b = "1"; c = "1"; d = "Tue Sep 22 22:56:15 CEST 2020"; e = "Tue Sep 22 22:56:36 CEST 2020"; f = "null" },
    @{ a = 129; b = "2"; c = "2"; d = "Tue Sep 22 22:54:38 CEST 2020"; e = "Tue Sep 22 22:56:54 CEST 2020"; f = "null" },
    @{ a = 130; b = "1"; c = "1"; d = "Tue Sep 22 22:56:36 CEST 2020"; e = "Tue Sep 22 22:57:20 CEST 2020"; f = "null" },
    @{ a = 131; b = "1"; c = "1"; d = "Tue Sep 22 22:57:21 CEST 2020"; e = "Tue Sep 22 22:57:23 CEST 2020"; f = "null" }
)

$startRow = 126
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row.a

    # B and C hold numeric-looking codes but must stay text (like the rest of
    # the sheet), so force text formatting, assign, then drop back to the
    # default style so no stray formatting remains on the cell.
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row.b
    $ws.Cells.Item($r, 2).Style = "Normal"

    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row.c
    $ws.Cells.Item($r, 3).Style = "Normal"

    $ws.Cells.Item($r, 4).Value = $row.d
    $ws.Cells.Item($r, 5).Value = $row.e
    $ws.Cells.Item($r, 6).Value = $row.f
}
